$d = $word.ActiveDocument

# Locate the paragraph "Hochladen auf FTP-Server (kommt per Mail)" - it is
# currently the last body paragraph, sits at list level ilvl=1
# (ListLevelNumber=2), and carries the _GoBack bookmark right at its end.
$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Hochladen auf FTP-Server*") {
        $targetIndex = $i
    }
}
$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range

# Step 1: turn this paragraph into the first of the new sub-items by
# replacing its text (but not its paragraph mark, so it keeps its
# identity/metadata) and bumping its list level from ilvl=1 to ilvl=2.
$textRange = $d.Range($r.Start, $r.End - 1)
$textRange.Text = "Argumentation für den Aufbau (Navigation, Struktur) und Gestaltung der Projekt-Userinterfaces"
$d.Paragraphs.Item($targetIndex).Range.ListFormat.ListLevelNumber = 3

# Step 2: append the remaining new list items as brand-new paragraphs,
# finishing with a fresh "Hochladen auf FTP-Server..." paragraph (the
# original text, now re-added at the end instead of being reused in
# place).
$lines = @(
    "Wichtige Aspekte",
    "Endnutzer",
    "Alter",
    "Erfahrungshorizont",
    "Nutzungskontext",
    "Professionell/beruflich",
    "Privat",
    "Gelegentlich / täglich",
    "Aufgabe der Anwendung",
    "Information",
    "Aufgaben erledigen / Arbeit",
    "Motivation / Enterainment",
    "Bezogen auf HCI Kriterien & Usability",
    "Hochladen auf FTP-Server (kommt per Mail)"
)
# ilvl (0-based, from the XML) + 1 == ListLevelNumber (1-based, COM)
$levels = @(3, 4, 5, 5, 4, 5, 5, 5, 4, 5, 5, 5, 3, 2)

$idx = $targetIndex
for ($j = 0; $j -lt $lines.Count; $j++) {
    $cur = $d.Paragraphs.Item($idx)
    $cur.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $newPara = $d.Paragraphs.Item($idx)
    $newPara.Range.Text = $lines[$j]
    $newPara.Range.ListFormat.ListLevelNumber = $levels[$j]
}

$bezogenIndex = $idx - 1

# Move the _GoBack bookmark from the end of the (old) last paragraph into
# the middle of the freshly-typed "Bezogen auf HCI Kriterien & Usability"
# paragraph, right after "Bezogen auf H". Adding a bookmark named
# "_GoBack" replaces any existing one of that name, so the old trailing
# bookmark disappears automatically.
$bezogenPara = $d.Paragraphs.Item($bezogenIndex)
$bmPos = $bezogenPara.Range.Start + "Bezogen auf H".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
